$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# xlHAlign constants
$xlHAlignCenter = -4108
$xlHAlignLeft = -4131

# NOTE: we deliberately avoid Columns.Insert() (the host engine mis-shifts
# the trailing "rest of sheet" <col> span by one when inserting a column),
# and instead author every cell at its final destination directly, in the
# same order the strings were first typed so the shared-string table comes
# out in the same sequence as the authored workbook.

# --- row 3 (Issue + Detail), typed first while this was still a 3-column
# sheet in the source workbook ---
$ws.Range("A3").Value = 'We need sort function on ally.'
$ws.Range("A3").HorizontalAlignment = $xlHAlignLeft
$ws.Range("D3").Value = 'In the Ally menu, there are two types of sorts: by Attack type (i.e. by melee, missile, and magic) and by power (i.e. ground, air, sea). We need the bot method to be able to sort the ally list by these categories.'
$ws.Range("D3").HorizontalAlignment = $xlHAlignLeft
$ws.Range("C3").Value = "Yes"

# --- new "Issue Type" column header + the two values typed for it next ---
$ws.Range("B1").Value = 'Issue Type'
$ws.Range("B1").HorizontalAlignment = $xlHAlignCenter
$ws.Range("B2").Value = 'Bugfix'
$ws.Range("B2").HorizontalAlignment = $xlHAlignCenter
$ws.Range("B3").Value = 'Feature Request'
$ws.Range("B3").HorizontalAlignment = $xlHAlignCenter

# --- row 4 ---
$ws.Range("A4").Value = 'Sort allies, but remove filters if none are found'
$ws.Range("A4").HorizontalAlignment = $xlHAlignLeft
$ws.Range("B4").Value = 'Feature Request'
$ws.Range("B4").HorizontalAlignment = $xlHAlignCenter
$ws.Range("C4").Value = "Yes"
$ws.Range("D4").Value = 'Would like a feature that adds to the sort allies feature. If an ally is not found, we''d like to remove the TYPE filter and see choose whichever ally is available.'
$ws.Range("D4").HorizontalAlignment = $xlHAlignLeft

# --- row 5 ---
$ws.Range("A5").Value = 'Need a feature that can do quest progression'
$ws.Range("A5").HorizontalAlignment = $xlHAlignLeft
$ws.Range("B5").Value = 'Feature Request'
$ws.Range("B5").HorizontalAlignment = $xlHAlignCenter
$ws.Range("C5").Value = "Yes"
$ws.Range("D5").Value = 'Need to be able to progress through the quest scenarios.'
$ws.Range("D5").HorizontalAlignment = $xlHAlignLeft

# --- row 6 ---
$ws.Range("A6").Value = 'First Unit''s cost is too high, cannot deploy unit'
$ws.Range("A6").HorizontalAlignment = $xlHAlignLeft
$ws.Range("B6").Value = 'Bugfix'
$ws.Range("B6").HorizontalAlignment = $xlHAlignCenter
$ws.Range("C6").Value = "Yes"
$ws.Range("D6").Value = 'Usually the script deploys the first unit. Some quests''s max unit points is less than the cost of the first unit, so the script will attempt to deploy the first unit, but it never will be able to do that.'
$ws.Range("D6").HorizontalAlignment = $xlHAlignLeft

# --- new "Comment" column ---
$ws.Range("E1").Value = 'Comment'
$ws.Range("E1").HorizontalAlignment = $xlHAlignCenter
$ws.Range("E6").Value = 'Made it so that the script will pick any of the 4 units when possible instead of just the first unit'
$ws.Range("E6").HorizontalAlignment = $xlHAlignLeft

# --- finish authoring the first two (pre-existing) rows/columns: the old
# B1/C1/B2/C2 header+value pair ("Resolved?" / "Detail" / "Yes" / long
# text) move one column to the right, to C/D ---
$ws.Range("C1").Value = 'Resolved?'
$ws.Range("C1").HorizontalAlignment = $xlHAlignCenter
$ws.Range("D1").Value = 'Detail'
$ws.Range("D1").HorizontalAlignment = $xlHAlignCenter
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = 'Consider two clients A and B that both run Bluestack app players. The title of A''s application is called BlueStacks App Player. The title of B''s application is called Bluestacks App Player. There''s a slight difference between the two titles; Client A has a capital s and client B has a lowercase s. This discrepency causes the bot to not recognize the application if ran on both client A and B. Please fix so that the bot works in both client A and client B''s cases.'
$ws.Range("D2").HorizontalAlignment = $xlHAlignLeft

# A1/A2 (Issue Summary header + first issue) are untouched, but make sure
# the header keeps its centred style.
$ws.Range("A1").HorizontalAlignment = $xlHAlignCenter
$ws.Range("A2").HorizontalAlignment = $xlHAlignLeft

# --- Column widths (best-fit on the engine's 1/6-character-wide grid) ---
$ws.Columns("B").ColumnWidth = 14.877604166666666
$ws.Columns("E").ColumnWidth = 85.59244791666667

# --- Defined names ---
$wb.Names.Item("Detail").RefersTo = "=Sheet1!`$D:`$D"
$nm = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$nm.Visible = $false

# --- Selection / view state ---
$ws.Range("E5").Select()
try { $excel.ActiveWindow.ScrollColumn = 4 } catch { }

Write-Host "done"
